$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "Errors:"
$ws.Range("A4").Value = "Row 1: Trùng Title. "
$ws.Range("A5").Value = "Row 2: Ngày bắt đầu không được là quá khứ. Định dạng ngày kết thúc không hợp lệ hoặc trống. Giá trị T không hợp lệ hoặc trống. "
